# Generate Report for Handback
#
# A new handback was produced for e2e/2e8e7e68-51be-4dbc-af07-02b407cd4542.md
# in both the zh-cn and de-de localization status tables (row 7), but the
# handback turned out to be stale (not built from the latest handoff), so
# the "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" / "Error Detail" columns are populated for that row on both
# locale sheets.

$wb = $excel.ActiveWorkbook

$latestHandbackUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5c99da7f9be14ae45be5bd43392d996632c6fa1b/e2e/2e8e7e68-51be-4dbc-af07-02b407cd4542.md"
$handbackDisplay = "2e8e7e68-51be-4dbc-af07-02b407cd4542.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b901adbf5a52e6b1cbff70f2febece0777f5e4c7/e2e/2e8e7e68-51be-4dbc-af07-02b407cd4542.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5c99da7f9be14ae45be5bd43392d996632c6fa1b/e2e/2e8e7e68-51be-4dbc-af07-02b407cd4542.md."

# --- zh-cn sheet, row 7 ---
$wsZh = $wb.Worksheets.Item("zh-cn")

# I7: Latest Target File -> becomes a hyperlink to the handback markdown file
$wsZh.Hyperlinks.Add($wsZh.Range("I7"), $latestHandbackUrl, "", "", $handbackDisplay)

# J7: Latest Handback File
$wsZh.Range("J7").Value = "2e8e7e68-51be-4dbc-af07-02b407cd4542.218d3db09e4c2256fb6346c4f5a827277eaa63c0.zh-cn.xlf"

# K7: Latest Handback DateTime
$wsZh.Range("K7").Value = "2016-09-06 07:04:58"

# P7: Error Detail
$wsZh.Range("P7").Value = $errorDetail

# --- de-de sheet, row 7 ---
$wsDe = $wb.Worksheets.Item("de-de")

# I7: Latest Target File -> becomes a hyperlink to the handback markdown file
$wsDe.Hyperlinks.Add($wsDe.Range("I7"), $latestHandbackUrl, "", "", $handbackDisplay)

# J7: Latest Handback File
$wsDe.Range("J7").Value = "2e8e7e68-51be-4dbc-af07-02b407cd4542.218d3db09e4c2256fb6346c4f5a827277eaa63c0.de-de.xlf"

# K7: Latest Handback DateTime
$wsDe.Range("K7").Value = "2016-09-06 07:05:15"

# P7: Error Detail
$wsDe.Range("P7").Value = $errorDetail
